# Daily attendance processing - 2026-01-08 11:35:28
# Swap the order of the "Recorded By" names in column G:
#   "System, dnasr281@gmail.com"  ->  "dnasr281@gmail.com, System"
# Only cells whose value is exactly the old combined string are touched;
# cells that only contain "System" or only "dnasr281@gmail.com" (and the
# header "Recorded By") are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
